$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 5
$ws.Range("H5").Value = 65.333336
$ws.Range("I5").Value = 65.333336
$ws.Range("K5").Value = 65.333336
$ws.Range("M5").Value = 49.666664
# row 12
$ws.Range("H12").Value = 426.5
$ws.Range("I12").Value = 389.75
$ws.Range("K12").Value = 389.75
$ws.Range("M12").Value = -219.75
# row 33
$ws.Range("H33").Value = 181.5
$ws.Range("I33").Value = 181.5
$ws.Range("K33").Value = 181.5
$ws.Range("M33").Value = 47.5
# row 38
$ws.Range("H38").Value = 882.4286
$ws.Range("I38").Value = 44.25
$ws.Range("K38").Value = 132.75
$ws.Range("M38").Value = 239.25
# row 53
$ws.Range("H53").Value = 175.85715
$ws.Range("I53").Value = 151.2
$ws.Range("K53").Value = 151.2
$ws.Range("M53").Value = 485.8
# row 86
$ws.Range("H86").Value = 3007.6667
$ws.Range("I86").Value = 1699
$ws.Range("K86").Value = 1699
$ws.Range("M86").Value = -576
# row 89
$ws.Range("H89").Value = 3007.6667
$ws.Range("I89").Value = 1699
$ws.Range("K89").Value = 8495
$ws.Range("M89").Value = -2879
# row 98
$ws.Range("H98").Value = 3130
$ws.Range("I98").Value = 3130
$ws.Range("K98").Value = 3130
$ws.Range("M98").Value = -1632
# row 122
$ws.Range("H122").Value = 3130
$ws.Range("I122").Value = 3130
$ws.Range("K122").Value = 9390
$ws.Range("M122").Value = -6940
# row 137
$ws.Range("H137").Value = 1766.091
$ws.Range("J137").Value = 1949.5
$ws.Range("L137").Value = 5848.5
$ws.Range("N137").Value = -10948.5
# row 138
$ws.Range("H138").Value = 3254.238
$ws.Range("J138").Value = 3999
$ws.Range("L138").Value = 11997
$ws.Range("N138").Value = -22277

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 43
$ws.Range("H43").Value = 40000
$ws.Range("I43").Value = 40000
$ws.Range("K43").Value = 40000
$ws.Range("M43").Value = -39687
# row 132
$ws.Range("H132").Value = 2456.8333
$ws.Range("I132").Value = 2748.4
$ws.Range("K132").Value = 8245.200000000001
$ws.Range("M132").Value = -5715.200000000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 2364.5
$ws.Range("I86").Value = 2364.5
$ws.Range("K86").Value = 2364.5
$ws.Range("M86").Value = -1241.5
# row 89
$ws.Range("H89").Value = 2364.5
$ws.Range("I89").Value = 2364.5
$ws.Range("K89").Value = 11822.5
$ws.Range("M89").Value = -6206.5
# row 94
$ws.Range("H94").Value = 1187.8
$ws.Range("I94").Value = 1187.8
$ws.Range("K94").Value = 1187.8
$ws.Range("M94").Value = -736.8
# row 95
$ws.Range("H95").Value = 21312
$ws.Range("J95").Value = 21312
$ws.Range("L95").Value = 21312
$ws.Range("N95").Value = -26804
# row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 3409.3333
$ws.Range("J62").Value = 3533.3333
$ws.Range("L62").Value = 3533.3333
$ws.Range("N62").Value = -4781.3333
# row 65
$ws.Range("H65").Value = 3409.3333
$ws.Range("J65").Value = 3533.3333
$ws.Range("L65").Value = 17666.6665
$ws.Range("N65").Value = -23906.6665
# row 107
$ws.Range("H107").Value = 671.25
$ws.Range("I107").Value = 671.25
$ws.Range("K107").Value = 671.25
$ws.Range("M107").Value = 1248.75
# row 122
$ws.Range("H122").Value = 7000
$ws.Range("J122").Value = 7000
$ws.Range("L122").Value = 21000
$ws.Range("N122").Value = -25900
# row 134
$ws.Range("H134").Value = 4736.75
$ws.Range("J134").Value = 5999
$ws.Range("L134").Value = 17997
$ws.Range("N134").Value = -23067
# row 140
$ws.Range("H140").Value = 65000
$ws.Range("J140").Value = 65000
$ws.Range("L140").Value = 65000
$ws.Range("N140").Value = -75360

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 2894
$ws.Range("I5").Value = 2894
$ws.Range("K5").Value = 8682
$ws.Range("M5").Value = -8570
# row 86
$ws.Range("H86").Value = 1300
$ws.Range("I86").Value = 600
$ws.Range("K86").Value = 1800
$ws.Range("M86").Value = -614
# row 89
$ws.Range("H89").Value = 1300
$ws.Range("I89").Value = 600
$ws.Range("K89").Value = 5400
$ws.Range("M89").Value = 528
# row 92
$ws.Range("H92").Value = 1412.25
$ws.Range("I92").Value = 2575
$ws.Range("J92").Value = 249.5
$ws.Range("K92").Value = 7725
$ws.Range("L92").Value = 748.5
$ws.Range("M92").Value = -6477
$ws.Range("N92").Value = -3244.5
# row 109
$ws.Range("H109").Value = 1113.5
$ws.Range("I109").Value = 227
$ws.Range("K109").Value = 681
$ws.Range("M109").Value = 359
# row 135
$ws.Range("H135").Value = 2894
$ws.Range("I135").Value = 2894
$ws.Range("K135").Value = 26046
$ws.Range("M135").Value = -23511
# row 139
$ws.Range("H139").Value = 5000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = $null
$ws.Range("N139").Value = -25280

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 584.5
$ws.Range("I102").Value = 584.5
$ws.Range("K102").Value = 584.5
$ws.Range("M102").Value = 1037.5
# row 122
$ws.Range("H122").Value = 3342.8333
$ws.Range("I122").Value = 3611.4
$ws.Range("K122").Value = 10834.2
$ws.Range("M122").Value = -8384.200000000001
# row 136
$ws.Range("H136").Value = 40659.8
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 40659.8
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 121979.4
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = -127079.4

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 6578.727
$ws.Range("I40").Value = 6420.75
$ws.Range("K40").Value = 6420.75
$ws.Range("M40").Value = -6284.75
# row 61
$ws.Range("H61").Value = 4500
$ws.Range("I61").Value = 4500
$ws.Range("K61").Value = 4500
$ws.Range("M61").Value = -4298
# row 101
$ws.Range("H101").Value = 13000
$ws.Range("J101").Value = 13000
$ws.Range("L101").Value = 13000
$ws.Range("N101").Value = -19490
# row 113
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330
# row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 12498
$ws.Range("I132").Value = 12498
$ws.Range("K132").Value = 37494
$ws.Range("M132").Value = -34964
